# Refresh crypto price/volume snapshot (GitHub Actions scrape, 2023-05-25 08:07:30 UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.332.25'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '''1.793.52'
$ws.Range('E3').Value = '  -1.92%  '
$ws.Range('D4').Value = '''1.007'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '''307.51'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('D7').Value = '''0.4531'
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('D8').Value = '''0.3595'
$ws.Range('E8').Value = '  -2.62%  '
$ws.Range('D9').Value = '''45.64'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').Value = '''0.07081'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = '''0.8858'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').Value = '''0.07827'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = '''19.43'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '''1.839.27'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '''5.291'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '''6.339'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').Value = '''84.76'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '''1.008'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Value = '''0.000008520'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').Value = '''14.26'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').Value = '''26.342.98'
$ws.Range('E22').Value = '  -2.03%  '
$ws.Range('D23').Value = '''4.986'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '''2.056.56'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D25').Value = '''10.50'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').Value = '''1.967'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '''152.18'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = '''17.84'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('E29').Value = '  +3.13%  '
$ws.Range('D30').Value = '''111.98'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D31').Value = '''4.865'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').Value = '''0.08692'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').Value = '''3.074'
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').Value = '''4.450'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '''2.721'
$ws.Range('E35').Value = '  +6.32%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.7210'
$ws.Range('E36').Value = '  -4.57%  '
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '''1.006'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''1.071'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.01927'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '''0.05103'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '''2.877'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = '''0.5077'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''6.855'
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '''0.1512'
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '''7.986'
$ws.Range('E46').Value = '  -3.54%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '''1.008'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = '''0.4628'
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''100.82'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''9.845'
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.578'
$ws.Range('E51').Value = '  -1.74%  '
